$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 55; this shifts existing rows 55-142 down to 56-143
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new record's data
$ws.Cells.Item(55, 1).Value = 3
$ws.Cells.Item(55, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(55, 3).Value = "Coquimbo"
$ws.Cells.Item(55, 4).Value = 44725
$ws.Cells.Item(55, 5).Value = 5
$ws.Cells.Item(55, 6).Value = 100112026
$ws.Cells.Item(55, 7).Value = "Haba"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 77
$ws.Cells.Item(55, 11).Value = 21000
$ws.Cells.Item(55, 12).Value = 22000
$ws.Cells.Item(55, 13).Value = 21506
$ws.Cells.Item(55, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(55, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(55, 16).Value = 860
$ws.Cells.Item(55, 17).Value = 25
$ws.Cells.Item(55, 18).Value = "Hortaliza"
